$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$overview.Range("G4").Value = "2016-09-06 17:20:47"
$dede.Range("H4").Value = "2016-09-06 17:20:47"

$zhcn.Range("H4").Value = "2016-09-06 17:20:41"
$zhcn.Range("K4").Value = "2016-09-06 17:21:01"

$dede.Range("K4").Value = "2016-09-06 17:21:23"
